$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 712; this shifts the existing rows 712-753
# down to 713-754 and extends the sheet dimension automatically.
$ws.Rows.Item(712).Insert()

# Fill in the newly inserted row 712 with the new entry.
# Force column A to be stored as literal text (matching the existing
# "yyyy/mm/dd" text entries below it) instead of being auto-converted
# to a date serial number by Excel's smart input parsing.
$ws.Range("A712").NumberFormat = "@"
$ws.Range("A712").Value = "2026/01/25"
$ws.Range("A712").Style = "Normal"

$ws.Range("B712").Value = "日"
$ws.Range("C712").Value = 19
$ws.Range("D712").Value = 201
